# Fill in the tester name for the single test-case row.
# S4 = 测试人员 (Tester) column — previously blank, now "黄烁".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("S4").Value = "黄烁"
